$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-17 01:21:19"

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
